# Auto: Update ETF Data
# Set "Share Change" (col G) and "Net Amount" (col H) to 0 for the
# data rows whose values changed in this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,24,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,46,47,48,49,51,52,53)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = 0
}
